$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7623068690299988
$ws.Range("B1").Value = 2.072154998779297
$ws.Range("C1").Value = 2.020151615142822
$ws.Range("D1").Value = 0.9882339239120483
$ws.Range("E1").Value = 1.127677202224731
